$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column H ("Save"), matching the formatting (bold, centered,
# thin border) of the other header cells in row 1 by copying from G1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" column values for rows 2-8
$saveValues = @(0, 1, 0, 0, 1, 1, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
